$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.009319804415525823
$ws.Range("D2").Value = 0.0107095245145814
$ws.Range("E2").Value = 0.4201940620647378
$ws.Range("F2").Value = 0.7392819203125214
$ws.Range("G2").Value = 0.594133829690378
$ws.Range("H2").Value = 0.6552490183734392
$ws.Range("I2").Value = 0.5535716855459185
$ws.Range("K2").Value = 1.384413480739227
$ws.Range("N2").Value = 0.9350266736026214

$ws.Range("C3").Value = 0.008256651033008211
$ws.Range("D3").Value = 0.01056964768328683
$ws.Range("E3").Value = 0.3663914740238141
$ws.Range("F3").Value = 0.7217621395818838
$ws.Range("G3").Value = 0.5760779366716093
$ws.Range("H3").Value = 0.6527258562016556
$ws.Range("I3").Value = 0.5411123819884267
$ws.Range("K3").Value = 1.215709099681476
$ws.Range("N3").Value = 0.945285378288645

$ws.Range("C4").Value = 0.007600189705975424
$ws.Range("D4").Value = 0.01048610810387629
$ws.Range("E4").Value = 0.3334698187378677
$ws.Range("F4").Value = 0.7116662727425194
$ws.Range("G4").Value = 0.5655967188800588
$ws.Range("H4").Value = 0.6516638329863156
$ws.Range("I4").Value = 0.5339640539766606
$ws.Range("K4").Value = 1.112131761093451
$ws.Range("N4").Value = 0.952040308515663

$ws.Range("C5").Value = 0.007331766573557275
$ws.Range("D5").Value = 0.01045265771041315
$ws.Range("E5").Value = 0.3200801055332647
$ws.Range("F5").Value = 0.7077172925392716
$ws.Range("G5").Value = 0.5614763219763432
$ws.Range("H5").Value = 0.6513530439956128
$ws.Range("I5").Value = 0.5311763373043945
$ws.Range("K5").Value = 1.069925126633052
$ws.Range("N5").Value = 0.9549076464699979

$ws.Range("C6").Value = 0.007287140635540368
$ws.Range("D6").Value = 0.01044713918771123
$ws.Range("E6").Value = 0.3178582629860358
$ws.Range("F6").Value = 0.7070715100739733
$ws.Range("G6").Value = 0.5608011996756801
$ws.Range("H6").Value = 0.6513087907862882
$ws.Range("I6").Value = 0.5307209803511128
$ws.Range("K6").Value = 1.062916867175772
$ws.Range("N6").Value = 0.9553906896581807

$ws.Range("C7").Value = 0.007596573319382571
$ws.Range("D7").Value = 0.0104856545772094
$ws.Range("E7").Value = 0.3332891376908833
$ws.Range("F7").Value = 0.7116123481621628
$ws.Range("G7").Value = 0.56554054108841
$ws.Range("H7").Value = 0.6516591483215848
$ws.Range("I7").Value = 0.5339259517480883
$ws.Range("K7").Value = 1.111562539644524
$ws.Range("N7").Value = 0.9520785142676615

$ws.Range("C8").Value = 0.008954002165111774
$ws.Range("D8").Value = 0.01066080973339467
$ws.Range("E8").Value = 0.401617953677416
$ws.Range("F8").Value = 0.7331031447846073
$ws.Range("G8").Value = 0.5877817074882046
$ws.Range("H8").Value = 0.6542775958384226
$ws.Range("I8").Value = 0.5491710459079044
$ws.Range("K8").Value = 1.326242525013356
$ws.Range("N8").Value = 0.9384692046118985

$ws.Range("C9").Value = 0.01158615017533293
$ws.Range("D9").Value = 0.01102282015004974
$ws.Range("E9").Value = 0.5366236924834311
$ws.Range("F9").Value = 0.7805486345121579
$ws.Range("G9").Value = 0.6362643134588239
$ws.Range("H9").Value = 0.6633029436360971
$ws.Range("I9").Value = 0.583089868278833
$ws.Range("K9").Value = 1.747314216361474
$ws.Range("N9").Value = 0.9154014513281581

$ws.Range("C10").Value = 0.0135012866743125
$ws.Range("D10").Value = 0.01130005651588561
$ws.Range("E10").Value = 0.6365971663988148
$ws.Range("F10").Value = 0.8187196815787843
$ws.Range("G10").Value = 0.6749479492888781
$ws.Range("H10").Value = 0.6723420488678471
$ws.Range("I10").Value = 0.6105247273710575
$ws.Range("K10").Value = 2.056791805384421
$ws.Range("N10").Value = 0.9006631086203498

$ws.Range("C11").Value = 0.0143683687485705
$ws.Range("D11").Value = 0.01142863265999594
$ws.Range("E11").Value = 0.682285124328402
$ws.Range("F11").Value = 0.8368218317667271
$ws.Range("G11").Value = 0.6932324467886701
$ws.Range("H11").Value = 0.6769851478801172
$ws.Range("I11").Value = 0.6235652940312093
$ws.Range("K11").Value = 2.197622782219639
$ws.Range("N11").Value = 0.8944388189284425

$ws.Range("C12").Value = 0.01469610493027318
$ws.Range("D12").Value = 0.01147767548582124
$ws.Range("E12").Value = 0.6996189600918115
$ws.Range("F12").Value = 0.8437841257211147
$ws.Range("G12").Value = 0.7002567573708518
$ws.Range("H12").Value = 0.6788203803855595
$ws.Range("I12").Value = 0.6285850442849465
$ws.Range("K12").Value = 2.250959396141013
$ws.Range("N12").Value = 0.8921509947156494

$ws.Range("C13").Value = 0.01462554841755548
$ws.Range("D13").Value = 0.01146709749087194
$ws.Range("E13").Value = 0.6958843068301519
$ws.Range("F13").Value = 0.8422798737697121
$ws.Range("G13").Value = 0.6987394574863117
$ws.Range("H13").Value = 0.6784216961947891
$ws.Range("I13").Value = 0.6275003073684644
$ws.Range("K13").Value = 2.239472090144147
$ws.Range("N13").Value = 0.8926406398258351

$ws.Range("C14").Value = 0.0143953440973732
$ws.Range("D14").Value = 0.01143266034903689
$ws.Range("E14").Value = 0.683710516687583
$ws.Range("F14").Value = 0.8373924639610948
$ws.Range("G14").Value = 0.6938083205252781
$ws.Range("H14").Value = 0.67713458701013
$ws.Range("I14").Value = 0.6239766313535995
$ws.Range("K14").Value = 2.202010677152998
$ws.Range("N14").Value = 0.8942492110822684

$ws.Range("C15").Value = 0.01425425752439935
$ws.Range("D15").Value = 0.0114116126966799
$ws.Range("E15").Value = 0.6762580698117233
$ws.Range("F15").Value = 0.8344128097046166
$ws.Range("G15").Value = 0.6908009742396928
$ws.Range("H15").Value = 0.6763562403123728
$ws.Range("I15").Value = 0.62182893076222
$ws.Range("K15").Value = 2.179065391569679
$ws.Range("N15").Value = 0.8952435203004896

$ws.Range("C16").Value = 0.01344453647347876
$ws.Range("D16").Value = 0.01129170325775775
$ws.Range("E16").Value = 0.633615794262866
$ws.Range("F16").Value = 0.8175516225931574
$ws.Range("G16").Value = 0.6737669741048364
$ws.Range("H16").Value = 0.6720493518679973
$ws.Range("I16").Value = 0.6096838557107844
$ws.Range("K16").Value = 2.047589142022503
$ws.Range("N16").Value = 0.9010795545090389

$ws.Range("C17").Value = 0.01294673127400614
$ws.Range("D17").Value = 0.01121877263192061
$ws.Range("E17").Value = 0.6075116897200417
$ws.Range("F17").Value = 0.8073977482704606
$ws.Range("G17").Value = 0.6634942607068979
$ws.Range("H17").Value = 0.6695436916320432
$ws.Range("I17").Value = 0.6023774622148892
$ws.Range("K17").Value = 1.966945069340909
$ws.Range("N17").Value = 0.9047828590131246

$ws.Range("C18").Value = 0.01266001974543229
$ws.Range("D18").Value = 0.01117705636314525
$ws.Range("E18").Value = 0.5925168586093008
$ws.Range("F18").Value = 0.8016268899659735
$ws.Range("G18").Value = 0.6576502482304818
$ws.Range("H18").Value = 0.6681524609199698
$ws.Range("I18").Value = 0.5982276968319411
$ws.Range("K18").Value = 1.920565249892547
$ws.Range("N18").Value = 0.9069580980743694

$ws.Range("C19").Value = 0.01256287815019164
$ws.Range("D19").Value = 0.01116297171814296
$ws.Range("E19").Value = 0.5874431451275655
$ws.Range("F19").Value = 0.7996848538675891
$ws.Range("G19").Value = 0.6556826104346953
$ws.Range("H19").Value = 0.6676899760745982
$ws.Range("I19").Value = 0.5968316749561211
$ws.Range("K19").Value = 1.904862614512979
$ws.Range("N19").Value = 0.9077023556802786

$ws.Range("C20").Value = 0.01299976367273814
$ws.Range("D20").Value = 0.0112265122639883
$ws.Range("E20").Value = 0.6102884726323055
$ws.Range("F20").Value = 0.8084714555303378
$ws.Range("G20").Value = 0.6645811134886799
$ws.Range("H20").Value = 0.6698052485173775
$ws.Range("I20").Value = 0.6031497810207469
$ws.Range("K20").Value = 1.975529300537119
$ws.Range("N20").Value = 0.9043839574785792

$ws.Range("C21").Value = 0.01446297734066349
$ws.Range("D21").Value = 0.01144276576541969
$ws.Range("E21").Value = 0.6872853423701741
$ws.Range("F21").Value = 0.8388250895056899
$ws.Range("G21").Value = 0.695253978504212
$ws.Range("H21").Value = 0.6775105477402121
$ws.Range("I21").Value = 0.6250093986026002
$ws.Range("K21").Value = 2.21301381208923
$ws.Range("N21").Value = 0.8937748562091201

$ws.Range("C22").Value = 0.01541571518939833
$ws.Range("D22").Value = 0.0115861627845355
$ws.Range("E22").Value = 0.7377996843248411
$ws.Range("F22").Value = 0.8592895742428652
$ws.Range("G22").Value = 0.7158862264073775
$ws.Range("H22").Value = 0.682995461403209
$ws.Range("I22").Value = 0.6397718491230648
$ws.Range("K22").Value = 2.36826440787496
$ws.Range("N22").Value = 0.887244487008239

$ws.Range("C23").Value = 0.01490755176374137
$ws.Range("D23").Value = 0.01150944016325184
$ws.Range("E23").Value = 0.710820729400524
$ws.Range("F23").Value = 0.8483095285648545
$ws.Range("G23").Value = 0.7048202909221288
$ws.Range("H23").Value = 0.6800267674097427
$ws.Range("I23").Value = 0.6318489716215652
$ws.Range("K23").Value = 2.285400483692626
$ws.Range("N23").Value = 0.890692924786471

$ws.Range("C24").Value = 0.01297578933387911
$ws.Range("D24").Value = 0.01122301251412949
$ws.Range("E24").Value = 0.6090330495418215
$ws.Range("F24").Value = 0.8079858246878899
$ws.Range("G24").Value = 0.6640895546295553
$ws.Range("H24").Value = 0.6696868451048203
$ws.Range("I24").Value = 0.602800457596004
$ws.Range("K24").Value = 1.971648420922634
$ws.Range("N24").Value = 0.9045641570818788

$ws.Range("C25").Value = 0.01087732962413668
$ws.Range("D25").Value = 0.01092290800962203
$ws.Range("E25").Value = 0.4999754519126185
$ws.Range("F25").Value = 0.7671374627383614
$ws.Range("G25").Value = 0.6226174862849803
$ws.Range("H25").Value = 0.6604410990626377
$ws.Range("I25").Value = 0.5734768340716272
$ws.Range("K25").Value = 1.633387543859669
$ws.Range("N25").Value = 0.9212541022304634
